$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the summary figures at the top of the statement.
# ---------------------------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 110978
# Cant. Trabajadores (worker count)
$ws.Range("C13").Value = 2
# Cant. Periodos (period count)
$ws.Range("F13").Value = 3

# ---------------------------------------------------------------------------
# 2. DAIRA MARIA NORIEGA TEHERAN's existing row (row 17) keeps the same
#    worker but the overdue period changes from 2506 to 2502.
# ---------------------------------------------------------------------------
$ws.Range("E17").Value = "2502"

# ---------------------------------------------------------------------------
# 3. The old row 18 (IVONNE CECILIA CONVERS ESPINOSA) and row 19 (CARLOS
#    ENRIQUE ZUREK CONVERS) are removed from the statement. In their place a
#    new final row is added repeating DAIRA MARIA NORIEGA TEHERAN for period
#    2501. First bring row 19's "last row" formatting onto row 18 (heavier
#    bottom border, no fill) since that's the style the final data row uses.
# ---------------------------------------------------------------------------
$fmtSrc = $ws.Range("B19:J19")
$fmtDst = $ws.Range("B18:J18")
$fmtSrc.Copy() | Out-Null
$fmtDst.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Now overwrite row 18's content with DAIRA's data for period 2501.
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047432398"
$ws.Range("D18").Value = "DAIRA MARIA NORIEGA TEHERAN"
$ws.Range("E18").Value = "2501"
$ws.Range("F18").Value = 20878
$ws.Range("G18").Value = 1423500
$ws.Range("H18").Value = ""
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""

# Delete the now-redundant row 19 completely; everything below shifts up
# (rows 24/25 -> 23/24) and the merged cell references follow automatically.
$ws.Rows("19").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. The shaded background fill that used to highlight the worker detail
#    rows (16-18) is removed, leaving those rows unfilled (still keeping
#    their borders/number formats).
# ---------------------------------------------------------------------------
$dataRows = $ws.Range("B16:G18")
$dataRows.Interior.Pattern = -4142   # xlPatternNone
$dataRows.Interior.PatternColorIndex = -4105  # xlAutomatic
